# ---------------------------------------------------------------------------
# B6-PowerPoint.pptx edit
#
# 1) Re-point the three "Income statement" tables (slides 14/15/16) from the
#    old local table style {34A60E2C-E633-4BC9-8CBC-EF1F356E82E7} to the new
#    style {F7D6B6BA-EB60-4411-A73D-29F6188563A2}. Table styles cannot be
#    assigned through the .Style property (read-only) -- PowerPoint requires
#    Table.ApplyStyle("{GUID}") for this.
#
# 2) Swap the presentation's theme colour palette: the "Integral" colour
#    scheme that the slide master currently uses becomes the stock "Office
#    Theme" colour scheme (the font scheme and format scheme are already
#    identical between the two theme parts in this deck, so only the 12
#    theme colours actually need to change).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------------
$oldStyle = "{34A60E2C-E633-4BC9-8CBC-EF1F356E82E7}"
$newStyle = "{F7D6B6BA-EB60-4411-A73D-29F6188563A2}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# --- 2) Theme colours -------------------------------------------------------
# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink. Target values are the stock "Office" colour scheme.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function ToRGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$tcs.Item(1).RGB  = ToRGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = ToRGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink
